# Generate Report for Handoff
# Update the "latest handoff"-related timestamp cells for the file
# e5662835-b3be-4055-aceb-755b4b53c24f across the Overview, zh-cn and
# de-de sheets (row 7 in each table), as produced by re-running the
# handback/handoff report generation.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-31 04:43:17"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-31 04:43:13"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-31 04:43:17"
